$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting existing data right
$ws.Columns.Item(1).Insert()

# New column A header/value
$ws.Range("A1").Value = "TabName"
$ws.Range("A2").Value = "CasesTab"

# Apply wrap text to C2 (the StatQuery cell, now shifted from B2 to C2)
$ws.Range("C2").WrapText = $true

# Column widths (values chosen so the engine's internal pixel rounding
# reproduces the target XML <col width> as closely as possible).
# Column 5 (was column 4 pre-insert) keeps its original width/bestFit
# untouched, so it is intentionally not set here.
$ws.Columns.Item(1).ColumnWidth = 10.0
$ws.Columns.Item(2).ColumnWidth = 75.0
$ws.Columns.Item(3).ColumnWidth = 127.66666666666667
$ws.Columns.Item(4).ColumnWidth = 69.5

# View settings
$excel.ActiveWindow.Zoom = 40
$ws.Range("C11").Select()
